# Update scraped_at timestamps (column K) on the "snapshot" sheet.
# Mirrors a re-scrape run ~4 hours after the prior snapshot (2025-11-19 03:xx -> 07:xx UTC).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("snapshot")

$ws.Range("K2").Value = "2025-11-19T07:03:27.845460+00:00"
$ws.Range("K3").Value = "2025-11-19T07:03:27.845498+00:00"
$ws.Range("K4").Value = "2025-11-19T07:03:27.845518+00:00"
$ws.Range("K5").Value = "2025-11-19T07:03:33.207622+00:00"
$ws.Range("K6").Value = "2025-11-19T07:03:33.207652+00:00"
$ws.Range("K7").Value = "2025-11-19T07:03:38.176178+00:00"
$ws.Range("K8").Value = "2025-11-19T07:03:43.721083+00:00"
$ws.Range("K9").Value = "2025-11-19T07:03:49.643389+00:00"
$ws.Range("K10").Value = "2025-11-19T07:03:49.643417+00:00"
$ws.Range("K11").Value = "2025-11-19T07:04:00.098071+00:00"
$ws.Range("K12").Value = "2025-11-19T07:04:05.448928+00:00"
$ws.Range("K13").Value = "2025-11-19T07:04:10.876331+00:00"
$ws.Range("K14").Value = "2025-11-19T07:04:10.876357+00:00"
$ws.Range("K15").Value = "2025-11-19T07:04:10.876374+00:00"
$ws.Range("K16").Value = "2025-11-19T07:04:15.853810+00:00"
$ws.Range("K17").Value = "2025-11-19T07:04:32.370178+00:00"
$ws.Range("K18").Value = "2025-11-19T07:04:37.828847+00:00"
$ws.Range("K19").Value = "2025-11-19T07:04:42.824636+00:00"
$ws.Range("K20").Value = "2025-11-19T07:04:42.824667+00:00"
$ws.Range("K21").Value = "2025-11-19T07:04:42.824686+00:00"
$ws.Range("K22").Value = "2025-11-19T07:04:42.824702+00:00"
$ws.Range("K23").Value = "2025-11-19T07:04:48.274853+00:00"
$ws.Range("K24").Value = "2025-11-19T07:04:48.274883+00:00"
$ws.Range("K25").Value = "2025-11-19T07:04:53.273312+00:00"
$ws.Range("K26").Value = "2025-11-19T07:04:53.273338+00:00"
$ws.Range("K27").Value = "2025-11-19T07:04:53.273357+00:00"
$ws.Range("K28").Value = "2025-11-19T07:04:58.744787+00:00"
$ws.Range("K29").Value = "2025-11-19T07:04:58.744814+00:00"
$ws.Range("K30").Value = "2025-11-19T07:05:03.852997+00:00"
$ws.Range("K31").Value = "2025-11-19T07:05:03.853024+00:00"
$ws.Range("K32").Value = "2025-11-19T07:05:03.853041+00:00"
$ws.Range("K33").Value = "2025-11-19T07:05:03.853057+00:00"
$ws.Range("K34").Value = "2025-11-19T07:05:03.853071+00:00"
$ws.Range("K35").Value = "2025-11-19T07:05:09.254832+00:00"
$ws.Range("K36").Value = "2025-11-19T07:05:09.254884+00:00"
$ws.Range("K37").Value = "2025-11-19T07:05:19.480005+00:00"
$ws.Range("K38").Value = "2025-11-19T07:05:19.480035+00:00"
$ws.Range("K39").Value = "2025-11-19T07:05:19.480058+00:00"
$ws.Range("K40").Value = "2025-11-19T07:05:24.561921+00:00"
